$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add data rows 49-52 ---
# Insert each new row by copying the row immediately above it (which carries
# the correct style + relative formulas), then overwrite the raw input
# values (A: date, B/C/D: daily counts). E/F/G/H/I stay as copied formulas
# and recalculate automatically.

$ws.Rows.Item(48).Copy() | Out-Null
$ws.Rows.Item(49).Insert(-4121) | Out-Null
$ws.Range("A49").Value = 43946
$ws.Range("B49").Value = 2
$ws.Range("C49").Value = 27
$ws.Range("D49").Value = 0

$ws.Rows.Item(49).Copy() | Out-Null
$ws.Rows.Item(50).Insert(-4121) | Out-Null
$ws.Range("A50").Value = 43947
$ws.Range("B50").Value = 4
$ws.Range("C50").Value = 22
$ws.Range("D50").Value = 1

$ws.Rows.Item(50).Copy() | Out-Null
$ws.Rows.Item(51).Insert(-4121) | Out-Null
$ws.Range("A51").Value = 43948
$ws.Range("B51").Value = 12
$ws.Range("C51").Value = 35
$ws.Range("D51").Value = 4

$ws.Rows.Item(51).Copy() | Out-Null
$ws.Rows.Item(52).Insert(-4121) | Out-Null
$ws.Range("A52").Value = 43949
$ws.Range("B52").Value = 7
$ws.Range("C52").Value = 35
$ws.Range("D52").Value = 2

# --- Add two trailing, essentially blank rows (53-54), matching the
# workbook author's formatted-but-empty rows at the bottom of the sheet ---

$ws.Range("A52").Copy() | Out-Null
$ws.Rows.Item(53).Insert(-4121) | Out-Null
$ws.Range("A53").ClearContents()

$ws.Range("A53").Copy() | Out-Null
$ws.Rows.Item(54).Insert(-4121) | Out-Null
$ws.Range("A54").ClearContents()

# --- Update the view so the sheet opens scrolled near the new bottom rows ---

$ws.Activate() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E60").Select() | Out-Null
